# Fix erroneous IFRS financial figures in the 나노메딕스 company_list sheet.
# Rows 2-6 (2014-2018 IFRS) had inflated values pasted in from the wrong
# scale/company; rows 7-9 (2019E-2021E forecast columns) were bogus actuals
# that should not have been populated at all, so their data cells are cleared.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (spreadsheet row 2): correct D:AJ figures
$ws.Cells.Item(2, 4).Value = 167
$ws.Cells.Item(2, 5).Value = 12
$ws.Cells.Item(2, 6).Value = 12
$ws.Cells.Item(2, 7).Value = 8
$ws.Cells.Item(2, 8).Value = 8
$ws.Cells.Item(2, 9).Value = 8
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 123
$ws.Cells.Item(2, 12).Value = 45
$ws.Cells.Item(2, 13).Value = 77
$ws.Cells.Item(2, 14).Value = 78
$ws.Cells.Item(2, 15).Value = -1
$ws.Cells.Item(2, 16).Value = 81
$ws.Cells.Item(2, 17).Value = 2
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = -3
$ws.Cells.Item(2, 20).Value = 0
$ws.Cells.Item(2, 21).Value = 2
$ws.Cells.Item(2, 22).Value = 8
$ws.Cells.Item(2, 23).Value = 7.43
$ws.Cells.Item(2, 24).Value = 4.5
$ws.Cells.Item(2, 25).Value = 10.79
$ws.Cells.Item(2, 26).Value = 6.82
$ws.Cells.Item(2, 27).Value = 58.77
$ws.Cells.Item(2, 28).Value = 96.22
$ws.Cells.Item(2, 29).Value = 50
$ws.Cells.Item(2, 30).Value = 16.73
$ws.Cells.Item(2, 31).Value = 484
$ws.Cells.Item(2, 32).Value = 1.73
$ws.Cells.Item(2, 33).Value = 0
$ws.Cells.Item(2, 34).Value = 0
$ws.Cells.Item(2, 35).Value = 0
$ws.Cells.Item(2, 36).Value = 16186294

# Row 3 (spreadsheet row 3): correct D:AJ figures
$ws.Cells.Item(3, 4).Value = 181
$ws.Cells.Item(3, 5).Value = 16
$ws.Cells.Item(3, 6).Value = 16
$ws.Cells.Item(3, 7).Value = 1
$ws.Cells.Item(3, 8).Value = 7
$ws.Cells.Item(3, 9).Value = 7
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 142
$ws.Cells.Item(3, 12).Value = 57
$ws.Cells.Item(3, 13).Value = 84
$ws.Cells.Item(3, 14).Value = 86
$ws.Cells.Item(3, 15).Value = -1
$ws.Cells.Item(3, 16).Value = 81
$ws.Cells.Item(3, 17).Value = 4
$ws.Cells.Item(3, 18).Value = 17
$ws.Cells.Item(3, 19).Value = 3
$ws.Cells.Item(3, 20).Value = 0
$ws.Cells.Item(3, 21).Value = 3
$ws.Cells.Item(3, 22).Value = 26
$ws.Cells.Item(3, 23).Value = 8.720000000000001
$ws.Cells.Item(3, 24).Value = 3.79
$ws.Cells.Item(3, 25).Value = 8.640000000000001
$ws.Cells.Item(3, 26).Value = 5.18
$ws.Cells.Item(3, 27).Value = 67.98999999999999
$ws.Cells.Item(3, 28).Value = 106.1
$ws.Cells.Item(3, 29).Value = 44
$ws.Cells.Item(3, 30).Value = 80.43000000000001
$ws.Cells.Item(3, 31).Value = 529
$ws.Cells.Item(3, 32).Value = 6.66
$ws.Cells.Item(3, 33).Value = 0
$ws.Cells.Item(3, 34).Value = 0
$ws.Cells.Item(3, 35).Value = 0
$ws.Cells.Item(3, 36).Value = 16186294

# Row 4 (spreadsheet row 4): correct D:AJ figures
$ws.Cells.Item(4, 4).Value = 384
$ws.Cells.Item(4, 5).Value = 43
$ws.Cells.Item(4, 6).Value = 41
$ws.Cells.Item(4, 7).Value = 27
$ws.Cells.Item(4, 8).Value = 21
$ws.Cells.Item(4, 9).Value = 21
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 218
$ws.Cells.Item(4, 12).Value = 80
$ws.Cells.Item(4, 13).Value = 138
$ws.Cells.Item(4, 14).Value = 140
$ws.Cells.Item(4, 15).Value = -1
$ws.Cells.Item(4, 16).Value = 88
$ws.Cells.Item(4, 17).Value = 19
$ws.Cells.Item(4, 18).Value = -19
$ws.Cells.Item(4, 19).Value = 30
$ws.Cells.Item(4, 20).Value = 3
$ws.Cells.Item(4, 21).Value = 16
$ws.Cells.Item(4, 22).Value = 21
$ws.Cells.Item(4, 23).Value = 11.09
$ws.Cells.Item(4, 24).Value = 5.4
$ws.Cells.Item(4, 25).Value = 18.37
$ws.Cells.Item(4, 26).Value = 11.5
$ws.Cells.Item(4, 27).Value = 57.74
$ws.Cells.Item(4, 28).Value = 149.24
$ws.Cells.Item(4, 29).Value = 125
$ws.Cells.Item(4, 30).Value = 21.13
$ws.Cells.Item(4, 31).Value = 790
$ws.Cells.Item(4, 32).Value = 3.34
$ws.Cells.Item(4, 33).Value = 0
$ws.Cells.Item(4, 34).Value = 0
$ws.Cells.Item(4, 35).Value = 0
$ws.Cells.Item(4, 36).Value = 17666353

# Row 5 (spreadsheet row 5): correct D:AJ figures
$ws.Cells.Item(5, 4).Value = 447
$ws.Cells.Item(5, 5).Value = -5
$ws.Cells.Item(5, 6).Value = -5
$ws.Cells.Item(5, 7).Value = -110
$ws.Cells.Item(5, 8).Value = -118
$ws.Cells.Item(5, 9).Value = -118
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 518
$ws.Cells.Item(5, 12).Value = 385
$ws.Cells.Item(5, 13).Value = 133
$ws.Cells.Item(5, 14).Value = 135
$ws.Cells.Item(5, 15).Value = -1
$ws.Cells.Item(5, 16).Value = 109
$ws.Cells.Item(5, 17).Value = -214
$ws.Cells.Item(5, 18).Value = -182
$ws.Cells.Item(5, 19).Value = 426
$ws.Cells.Item(5, 20).Value = 114
$ws.Cells.Item(5, 21).Value = -328
$ws.Cells.Item(5, 22).Value = 274
$ws.Cells.Item(5, 23).Value = -1.05
$ws.Cells.Item(5, 24).Value = -26.34
$ws.Cells.Item(5, 25).Value = -85.81999999999999
$ws.Cells.Item(5, 26).Value = -31.95
$ws.Cells.Item(5, 27).Value = 288.37
$ws.Cells.Item(5, 28).Value = 105.55
$ws.Cells.Item(5, 29).Value = -579
$ws.Cells.Item(5, 30).Value = -10.79
$ws.Cells.Item(5, 31).Value = 615
$ws.Cells.Item(5, 32).Value = 10.17
$ws.Cells.Item(5, 33).Value = 0
$ws.Cells.Item(5, 34).Value = 0
$ws.Cells.Item(5, 35).Value = 0
$ws.Cells.Item(5, 36).Value = 21886139

# Row 6 (spreadsheet row 6): correct D:AJ figures
$ws.Cells.Item(6, 4).Value = 398
$ws.Cells.Item(6, 5).Value = -114
$ws.Cells.Item(6, 6).Value = -114
$ws.Cells.Item(6, 7).Value = -83
$ws.Cells.Item(6, 8).Value = -83
$ws.Cells.Item(6, 9).Value = -67
$ws.Cells.Item(6, 11).Value = 913
$ws.Cells.Item(6, 12).Value = 528
$ws.Cells.Item(6, 13).Value = 385
$ws.Cells.Item(6, 14).Value = 325
$ws.Cells.Item(6, 16).Value = 139
$ws.Cells.Item(6, 17).Value = -65
$ws.Cells.Item(6, 18).Value = -292
$ws.Cells.Item(6, 19).Value = 356
$ws.Cells.Item(6, 20).Value = 75
$ws.Cells.Item(6, 21).Value = -139
$ws.Cells.Item(6, 22).Value = 401
$ws.Cells.Item(6, 23).Value = -28.77
$ws.Cells.Item(6, 24).Value = -20.96
$ws.Cells.Item(6, 25).Value = -29.34
$ws.Cells.Item(6, 26).Value = -11.65
$ws.Cells.Item(6, 27).Value = 137.16
$ws.Cells.Item(6, 28).Value = 175.47
$ws.Cells.Item(6, 29).Value = -280
$ws.Cells.Item(6, 30).Value = -20.19
$ws.Cells.Item(6, 31).Value = 1170
$ws.Cells.Item(6, 32).Value = 4.84
$ws.Cells.Item(6, 33).Value = 0
$ws.Cells.Item(6, 34).Value = 0
$ws.Cells.Item(6, 35).Value = 0
$ws.Cells.Item(6, 36).Value = 27805125

# Row 7: remove the bogus financial data, keep only the index/label columns
$ws.Range("D7:AJ7").ClearContents()

# Row 8: remove the bogus financial data, keep only the index/label columns
$ws.Range("D8:AJ8").ClearContents()

# Row 9: remove the bogus financial data, keep only the index/label columns
$ws.Range("D9:AJ9").ClearContents()
